$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = "`n"
$content = "COFFEE MACHIEN" + $nl + "Maybe a coffee machien would help in productivty" + $nl + $nl + "By : Abdullah Elrouby" + $nl + "ID : HRAE1706" + $nl + "Date : 05-Jul-24 3:28:42 AM"

$ws.Range("A4").Value = "COFFEE MACHIEN"
$ws.Range("B4").Value = "Review"
$ws.Range("C4").Value = $content
$ws.Range("D4").Value = "Accepted"
$ws.Range("E4").Value = "05-Jul-24 3:28:42 AM"

$v = $ws.Range("C4").Value2
Write-Host "C4 value2: [$v]"
